$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 45; this pushes existing rows 45-52 down to 46-53.
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new data record.
$ws.Cells.Item(45, 1).Value = 11
$ws.Cells.Item(45, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(45, 3).Value = "Bíobío"
$ws.Cells.Item(45, 4).Value = 44946
$ws.Cells.Item(45, 4).NumberFormat = $ws.Cells.Item(46, 4).NumberFormat
$ws.Cells.Item(45, 5).Value = 8
$ws.Cells.Item(45, 6).Value = 100112030
$ws.Cells.Item(45, 7).Value = "Poroto granado"
$ws.Cells.Item(45, 8).Value = "Sin especificar"
$ws.Cells.Item(45, 9).Value = "Primera"
$ws.Cells.Item(45, 10).Value = 250
$ws.Cells.Item(45, 11).Value = 28000
$ws.Cells.Item(45, 12).Value = 30000
$ws.Cells.Item(45, 13).Value = 28800
$ws.Cells.Item(45, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(45, 15).Value = "Región del Maule"
$ws.Cells.Item(45, 16).Value = 1152
$ws.Cells.Item(45, 17).Value = 25
$ws.Cells.Item(45, 18).Value = "Hortaliza"
